# Apply the cryptos-list refresh: updated prices (col D) and
# 1h volume deltas (col E) for each row, plus the row 28/29 and
# 32/33 coin swaps (PancakeSwap<->RenderToken, Kaspa<->Fetch.AI).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.248.51'
$ws.Range("E2").Value = '  +1.19%  '

$ws.Range("D3").Value = '3.563.30'
$ws.Range("E3").Value = '  +4.88%  '

$ws.Range("E4").Value = '  -0.01%  '

$c = $ws.Range("D5")
$c.Value = '''606.61'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.04%  '

$c = $ws.Range("D6")
$c.Value = '''145.02'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +2.91%  '

$ws.Range("D7").Value = '3.561.64'
$ws.Range("E7").Value = '  +4.84%  '

$ws.Range("E8").Value = '  +0.27%  '

$ws.Range("E9").Value = '  +3.66%  '

$ws.Range("E10").Value = '  +1.79%  '

$ws.Range("E11").Value = '  +1.94%  '

$ws.Range("E12").Value = '  +1.89%  '

$ws.Range("D13").Value = '4.164.83'
$ws.Range("E13").Value = '  +4.72%  '

$ws.Range("E14").Value = '  +4.17%  '

$c = $ws.Range("D15")
$c.Value = '''30.12'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.55%  '

$ws.Range("D16").Value = '3.558.15'
$ws.Range("E16").Value = '  +4.77%  '

$ws.Range("D17").Value = '66.325.17'
$ws.Range("E17").Value = '  +1.33%  '

$ws.Range("E18").Value = '  -0.81%  '

$c = $ws.Range("D19")
$c.Value = '''11.51'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +10.33%  '

$c = $ws.Range("D20")
$c.Value = '''6.22'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +1.73%  '

$c = $ws.Range("D21")
$c.Value = '''14.99'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.88%  '

$c = $ws.Range("D22")
$c.Value = '''431.91'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +3.74%  '

$c = $ws.Range("D23")
$c.Value = '''0.610'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +5.33%  '

$c = $ws.Range("D24")
$c.Value = '''78.83'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.89%  '

$ws.Range("D25").Value = '3.703.83'
$ws.Range("E25").Value = '  +4.65%  '

$ws.Range("E26").Value = '  -0.05%  '

$ws.Range("E27").Value = '  +7.69%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D28")
$c.Value = '''8.04'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +2.57%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D29")
$c.Value = '''2.51'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +4.31%  '

$c = $ws.Range("D30")
$c.Value = '''9.20'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.15%  '

$c = $ws.Range("D31")
$c.Value = '''0.999'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.04%  '

$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D32")
$c.Value = '''1.48'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.35%  '

$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D33")
$c.Value = '''0.160'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.34%  '

$ws.Range("D34").Value = '3.554.47'
$ws.Range("E34").Value = '  +4.69%  '

$c = $ws.Range("D35")
$c.Value = '''25.44'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +4.05%  '

$ws.Range("E36").Value = '  +0.03%  '

$c = $ws.Range("D37")
$c.Value = '''1.75'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +3.04%  '

$ws.Range("E38").Value = '  +4.77%  '

$c = $ws.Range("D39")
$c.Value = '''5.65'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.37%  '

$c = $ws.Range("D40")
$c.Value = '''0.999'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.02%  '

$c = $ws.Range("D41")
$c.Value = '''171.56'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.40%  '

$c = $ws.Range("D42")
$c.Value = '''0.0854'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.64%  '

$c = $ws.Range("D43")
$c.Value = '''5.21'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +3.03%  '

$c = $ws.Range("D44")
$c.Value = '''0.898'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +3.26%  '

$c = $ws.Range("D45")
$c.Value = '''1.96'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +2.50%  '

$c = $ws.Range("D46")
$c.Value = '''46.10'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.35%  '

$ws.Range("E47").Value = '  +4.57%  '

$c = $ws.Range("D48")
$c.Value = '''25.92'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.29%  '

$c = $ws.Range("D49")
$c.Value = '''2.38'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +3.90%  '

$ws.Range("E50").Value = '  +0.78%  '

$c = $ws.Range("D51")
$c.Value = '''0.954'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +4.05%  '
